# Regenerate save_data column G ("K") values to use actual strikeouts (K)
# instead of the prior "Strike#" derived values, per recalculated
# std/mean stats and s_vals write-back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 2
    11 = 0
    13 = 1
    16 = 1
    17 = 1
    18 = 2
    19 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
